$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Row 12 (B12:E12): "Fatima" / 19 / "Femenino" -> "Fabian" / 19 / "masculino" / Universidad
$ws.Range("B12").Value = "Fabian"
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = "masculino"
$ws.Range("E12").Value = "Universidad"

$ws.Activate()
$ws.Range("D12").Select()
$ws.Application.ActiveWindow.Zoom = 286
